$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("statistics")

# New header cells - text first, then copy formatting from the existing
# header style (bold, bordered, centered) so the same style index is reused.
$ws.Range("D1").Value = "Yes"
$ws.Range("E1").Value = "No"

$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: D = "Qualiperf acknowledged" count (Yes), E = remainder (No)
$data = @(
    @(22, 11),
    @(7, 2),
    @(1, 1),
    @(0, 5),
    @(19, 2),
    @(26, 7),
    @(4, 22)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $data[$i][0]
    $ws.Cells.Item($row, 5).Value = $data[$i][1]
}
